# "Hardware/Purchase list.xlsx" update
# - Row 23 (SparkFun #8506, "2X5 shrouded header"): quantity 10 -> 5
# - New row 30: SWD (2x5 1.27mm) Cable Breakout Board, qty 5 @ $1.95, with
#   a hyperlink to the Adafruit product page
# - Totals recalculate automatically

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23: quantity changed from 10 to 5 ---------------------------------
$ws.Range("D23").Value = 5

# --- New row 30 --------------------------------------------------------------
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "SWD (2x5 1.27mm) Cable Breakout Board"
$ws.Range("D30").Value = 5
$ws.Range("E30").Value = 1.95
$ws.Range("F30").Formula = "=E30*D30"

$url = "https://www.adafruit.com/product/2743?hidden=yes&main_page=product_info&products_id=2743"
$ws.Range("G30").Value = $url
$ws.Hyperlinks.Add($ws.Range("G30"), $url)
# match the look of the other hyperlink cells in column G
$ws.Range("G30").Style = $ws.Range("G29").Style

# --- recalc + reflect the author's final cursor / view position ------------
$excel.Calculate()
$ws.Range("D24").Select() | Out-Null
